# Auto-generated edit script: updates Leve profit/price data cells
# across multiple sheets (ALC, ARM, BSM, CUL, GSM, LTW, WVR) to match
# refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 480.4
$ws.Range("I2").Value = 377.36365
$ws.Range("K2").Value = 377.36365
$ws.Range("M2").Value = -264.36365
$ws.Range("H19").Value = 1756.6
$ws.Range("I19").Value = 1338.25
$ws.Range("K19").Value = 1338.25
$ws.Range("M19").Value = -1163.25
$ws.Range("H62").Value = 5617.25
$ws.Range("I62").Value = 4267
$ws.Range("K62").Value = 4267
$ws.Range("M62").Value = -3643
$ws.Range("H65").Value = 5617.25
$ws.Range("I65").Value = 4267
$ws.Range("K65").Value = 21335
$ws.Range("M65").Value = -18215
$ws.Range("H98").Value = 521.4375
$ws.Range("I98").Value = 422.93332
$ws.Range("J98").Value = 1999
$ws.Range("K98").Value = 422.93332
$ws.Range("L98").Value = 1999
$ws.Range("M98").Value = 1075.06668
$ws.Range("N98").Value = -4995
$ws.Range("H122").Value = 521.4375
$ws.Range("I122").Value = 422.93332
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 1268.79996
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = 1181.20004
$ws.Range("N122").Value = -10897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1876
$ws.Range("I14").Value = 2101.6667
$ws.Range("J14").Value = 1199
$ws.Range("K14").Value = 2101.6667
$ws.Range("L14").Value = 1199
$ws.Range("M14").Value = -1926.6667
$ws.Range("N14").Value = -1549
$ws.Range("H16").Value = 286303.44
$ws.Range("I16").Value = 400561.2
$ws.Range("J16").Value = 659
$ws.Range("K16").Value = 400561.2
$ws.Range("L16").Value = 659
$ws.Range("M16").Value = -400274.2
$ws.Range("N16").Value = -1233
$ws.Range("H93").Value = 18663.334
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 18663.334
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 18663.334
$ws.Range("N93").Value = -23655.334
$ws.Range("M93").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 1017.3333
$ws.Range("J19").Value = 1017.3333
$ws.Range("L19").Value = 1017.3333
$ws.Range("N19").Value = -1363.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1933.3334
$ws.Range("I51").Value = 400
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 1200
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = -740
$ws.Range("N51").Value = -15920
$ws.Range("H68").Value = 2299.4
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2299.4
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6898.200000000001
$ws.Range("N68").Value = -8520.200000000001
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 2299.4
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2299.4
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 20694.6
$ws.Range("N71").Value = -28806.6
$ws.Range("M71").ClearContents()
$ws.Range("H75").Value = 2533.6667
$ws.Range("I75").Value = 2302.1667
$ws.Range("J75").Value = 2765.1667
$ws.Range("K75").Value = 6906.500100000001
$ws.Range("L75").Value = 8295.500100000001
$ws.Range("M75").Value = -5908.500100000001
$ws.Range("N75").Value = -10291.5001
$ws.Range("H78").Value = 2533.6667
$ws.Range("I78").Value = 2302.1667
$ws.Range("J78").Value = 2765.1667
$ws.Range("K78").Value = 20719.5003
$ws.Range("L78").Value = 24886.5003
$ws.Range("M78").Value = -15727.5003
$ws.Range("N78").Value = -34870.5003
$ws.Range("H80").Value = 6000
$ws.Range("J80").Value = 6000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -19872
$ws.Range("H83").Value = 6000
$ws.Range("J83").Value = 6000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -63360
$ws.Range("H92").Value = 511.33334
$ws.Range("I92").Value = 623
$ws.Range("J92").Value = 399.66666
$ws.Range("K92").Value = 1869
$ws.Range("L92").Value = 1198.99998
$ws.Range("M92").Value = -621
$ws.Range("N92").Value = -3694.99998
$ws.Range("H121").Value = 1131.0625
$ws.Range("I121").Value = 682.5
$ws.Range("J121").Value = 1579.625
$ws.Range("K121").Value = 2047.5
$ws.Range("L121").Value = 4738.875
$ws.Range("M121").Value = -737.5
$ws.Range("N121").Value = -7358.875
$ws.Range("H129").Value = 3052.9395
$ws.Range("I129").Value = 3214.8333
$ws.Range("K129").Value = 9644.499899999999
$ws.Range("M129").Value = -4644.499899999999
$ws.Range("H134").Value = 4210.75
$ws.Range("J134").Value = 14999.333
$ws.Range("L134").Value = 44997.999
$ws.Range("N134").Value = -55137.999
$ws.Range("H140").Value = 1910.75
$ws.Range("I140").Value = 1041.6
$ws.Range("K140").Value = 3124.8
$ws.Range("M140").Value = 2055.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4423.364
$ws.Range("I70").Value = 4295.4443
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 4295.4443
$ws.Range("L70").Value = 4999
$ws.Range("M70").Value = -4025.4443
$ws.Range("N70").Value = -5539
$ws.Range("H73").Value = 4423.364
$ws.Range("I73").Value = 4295.4443
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 4295.4443
$ws.Range("L73").Value = 4999
$ws.Range("M73").Value = -3359.4443
$ws.Range("N73").Value = -6871
$ws.Range("H132").Value = 2999.3547
$ws.Range("I132").Value = 2499.7368
$ws.Range("J132").Value = 3790.4167
$ws.Range("K132").Value = 7499.2104
$ws.Range("L132").Value = 11371.2501
$ws.Range("M132").Value = -4969.2104
$ws.Range("N132").Value = -16431.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 382.3
$ws.Range("I16").Value = 380.33334
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 380.33334
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -210.33334
$ws.Range("N16").Value = -740
$ws.Range("H61").Value = 5166
$ws.Range("I61").Value = 4441.8335
$ws.Range("J61").Value = 6252.25
$ws.Range("K61").Value = 4441.8335
$ws.Range("L61").Value = 6252.25
$ws.Range("M61").Value = -4239.8335
$ws.Range("N61").Value = -6656.25
$ws.Range("H113").Value = 5166
$ws.Range("I113").Value = 4441.8335
$ws.Range("J113").Value = 6252.25
$ws.Range("K113").Value = 4441.8335
$ws.Range("L113").Value = 6252.25
$ws.Range("M113").Value = -2271.8335
$ws.Range("N113").Value = -10592.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 14328184
$ws.Range("I86").Value = 25036824
$ws.Range("J86").Value = 49998.668
$ws.Range("K86").Value = 25036824
$ws.Range("L86").Value = 49998.668
$ws.Range("M86").Value = -25035701
$ws.Range("N86").Value = -52244.668
$ws.Range("H89").Value = 14328184
$ws.Range("I89").Value = 25036824
$ws.Range("J89").Value = 49998.668
$ws.Range("K89").Value = 125184120
$ws.Range("L89").Value = 249993.34
$ws.Range("M89").Value = -125178504
$ws.Range("N89").Value = -261225.34

Write-Host "Done updating Leve profit data."
